$d = $word.ActiveDocument

# Merge "Versi" + "on" runs into a single "Version" run.
# (Setting the same text is a no-op for Word's run-merging, so first
# change it to a placeholder, then correct it, forcing a real merge.)
$rVersion = $d.Range(0, 7)
$rVersion.Text = "Versio_"
$rVersion = $d.Range(0, 7)
$rVersion.Text = "Version"

# Change " 2" to " 1." within its own run (keeps proofErr/bookmark intact).
$rNum = $d.Range(8, 9)
$rNum.Text = "1."

# Remove the trailing "." run that followed the bookmark.
$rDot = $d.Range(10, 11)
$rDot.Delete()
